$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 998
$ws.Range("J19").Value = 998
$ws.Range("L19").Value = 998
$ws.Range("N19").Value = -1348
$ws.Range("H70").Value = 1598.2667
$ws.Range("I70").Value = 891.6667
$ws.Range("J70").Value = 1774.9166
$ws.Range("K70").Value = 2675.0001
$ws.Range("L70").Value = 5324.7498
$ws.Range("M70").Value = -2405.0001
$ws.Range("N70").Value = -5864.7498
$ws.Range("H73").Value = 1598.2667
$ws.Range("I73").Value = 891.6667
$ws.Range("J73").Value = 1774.9166
$ws.Range("K73").Value = 2675.0001
$ws.Range("L73").Value = 5324.7498
$ws.Range("M73").Value = -1739.0001
$ws.Range("N73").Value = -7196.7498
$ws.Range("H76").Value = 3967
$ws.Range("J76").Value = 3950
$ws.Range("L76").Value = 3950
$ws.Range("N76").Value = -4580
$ws.Range("H79").Value = 3967
$ws.Range("J79").Value = 3950
$ws.Range("L79").Value = 3950
$ws.Range("N79").Value = -6134
$ws.Range("H111").Value = 403.7
$ws.Range("I111").Value = 422.77777
$ws.Range("K111").Value = 1268.33331
$ws.Range("M111").Value = 1798.66669
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = -11508
$ws.Range("H137").Value = 1702.7222
$ws.Range("I137").Value = 1567.5714
$ws.Range("K137").Value = 4702.7142
$ws.Range("M137").Value = -2152.7142
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5089.8887
$ws.Range("I63").Value = 829.1429000000001
$ws.Range("K63").Value = 829.1429000000001
$ws.Range("M63").Value = -143.1429000000001
$ws.Range("H66").Value = 5089.8887
$ws.Range("I66").Value = 829.1429000000001
$ws.Range("K66").Value = 4145.7145
$ws.Range("M66").Value = -713.7145
$ws.Range("H74").Value = 18623.584
$ws.Range("I74").Value = 20184
$ws.Range("J74").Value = 15502.75
$ws.Range("K74").Value = 20184
$ws.Range("L74").Value = 15502.75
$ws.Range("M74").Value = -19310
$ws.Range("N74").Value = -17250.75
$ws.Range("H77").Value = 18623.584
$ws.Range("I77").Value = 20184
$ws.Range("J77").Value = 15502.75
$ws.Range("K77").Value = 100920
$ws.Range("L77").Value = 77513.75
$ws.Range("M77").Value = -96552
$ws.Range("N77").Value = -86249.75
$ws.Range("H97").Value = 1753.1578
$ws.Range("I97").Value = 449.15384
$ws.Range("K97").Value = 449.15384
$ws.Range("M97").Value = 46.84616
$ws.Range("H110").Value = 2316.389
$ws.Range("I110").Value = 980.9375
$ws.Range("K110").Value = 980.9375
$ws.Range("M110").Value = 1064.0625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 60074
$ws.Range("J35").Value = 60074
$ws.Range("L35").Value = 60074
$ws.Range("N35").Value = -60694
$ws.Range("H82").Value = 27806.5
$ws.Range("I82").Value = 3614.25
$ws.Range("J82").Value = 51998.75
$ws.Range("K82").Value = 3614.25
$ws.Range("L82").Value = 51998.75
$ws.Range("M82").Value = -3231.25
$ws.Range("N82").Value = -52764.75
$ws.Range("H85").Value = 27806.5
$ws.Range("I85").Value = 3614.25
$ws.Range("J85").Value = 51998.75
$ws.Range("K85").Value = 3614.25
$ws.Range("L85").Value = 51998.75
$ws.Range("M85").Value = -2288.25
$ws.Range("N85").Value = -54650.75
$ws.Range("H86").Value = 2991.4285
$ws.Range("I86").Value = 1748
$ws.Range("K86").Value = 1748
$ws.Range("M86").Value = -625
$ws.Range("H89").Value = 2991.4285
$ws.Range("I89").Value = 1748
$ws.Range("K89").Value = 8740
$ws.Range("M89").Value = -3124
$ws.Range("H134").Value = 3039.7
$ws.Range("I134").Value = 3039.7
$ws.Range("K134").Value = 9119.099999999999
$ws.Range("M134").Value = -6584.099999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5178.8
$ws.Range("I58").Value = 2347.5
$ws.Range("J58").Value = 7066.3335
$ws.Range("K58").Value = 2347.5
$ws.Range("L58").Value = 7066.3335
$ws.Range("M58").Value = -2144.5
$ws.Range("N58").Value = -7472.3335
$ws.Range("H62").Value = 9666.333000000001
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 9499.5
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 9499.5
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -10747.5
$ws.Range("H65").Value = 9666.333000000001
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 9499.5
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 47497.5
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -53737.5
$ws.Range("H86").Value = 3913.8125
$ws.Range("I86").Value = 4868.143
$ws.Range("J86").Value = 3171.5557
$ws.Range("K86").Value = 4868.143
$ws.Range("L86").Value = 3171.5557
$ws.Range("M86").Value = -3745.143
$ws.Range("N86").Value = -5417.5557
$ws.Range("H89").Value = 3913.8125
$ws.Range("I89").Value = 4868.143
$ws.Range("J89").Value = 3171.5557
$ws.Range("K89").Value = 24340.715
$ws.Range("L89").Value = 15857.7785
$ws.Range("M89").Value = -18724.715
$ws.Range("N89").Value = -27089.7785
$ws.Range("H122").Value = 6598.8
$ws.Range("I122").Value = 5248.5
$ws.Range("J122").Value = 12000
$ws.Range("K122").Value = 15745.5
$ws.Range("L122").Value = 36000
$ws.Range("M122").Value = -13295.5
$ws.Range("N122").Value = -40900
$ws.Range("H132").Value = 4939.4
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4939.4
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 14818.2
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -19878.2
$ws.Range("H136").Value = 5178.8
$ws.Range("I136").Value = 2347.5
$ws.Range("J136").Value = 7066.3335
$ws.Range("K136").Value = 7042.5
$ws.Range("L136").Value = 21199.0005
$ws.Range("M136").Value = -4492.5
$ws.Range("N136").Value = -26299.0005
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4158.2
$ws.Range("I80").Value = 3597.3333
$ws.Range("J80").Value = 4999.5
$ws.Range("K80").Value = 3597.3333
$ws.Range("L80").Value = 4999.5
$ws.Range("M80").Value = -2599.3333
$ws.Range("N80").Value = -6995.5
$ws.Range("H83").Value = 4158.2
$ws.Range("I83").Value = 3597.3333
$ws.Range("J83").Value = 4999.5
$ws.Range("K83").Value = 17986.6665
$ws.Range("L83").Value = 24997.5
$ws.Range("M83").Value = -12994.6665
$ws.Range("N83").Value = -34981.5
$ws.Range("H113").Value = 1796.8182
$ws.Range("I113").Value = 1594
$ws.Range("J113").Value = 2337.6667
$ws.Range("K113").Value = 1594
$ws.Range("L113").Value = 2337.6667
$ws.Range("M113").Value = 576
$ws.Range("N113").Value = -6677.6667
$ws.Range("H132").Value = 4007.2
$ws.Range("I132").Value = 3515.6
$ws.Range("K132").Value = 10546.8
$ws.Range("M132").Value = -8016.799999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 595.3333
$ws.Range("I55").Value = 544.75
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 544.75
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = -371.75
$ws.Range("N55").Value = -1346
$ws.Range("H68").Value = 2669.4443
$ws.Range("I68").Value = 2669.4443
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2669.4443
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1920.4443
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 2669.4443
$ws.Range("I71").Value = 2669.4443
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 13347.2215
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -9603.2215
$ws.Range("N71").Value = ""
$ws.Range("H132").Value = 4499.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = -18558.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 597.9
$ws.Range("I100").Value = 590.1429000000001
$ws.Range("J100").Value = 616
$ws.Range("K100").Value = 1180.2858
$ws.Range("L100").Value = 1232
$ws.Range("M100").Value = -639.2858000000001
$ws.Range("N100").Value = -2314
$ws.Range("H107").Value = 534.1429000000001
$ws.Range("I107").Value = 472.125
$ws.Range("K107").Value = 1416.375
$ws.Range("M107").Value = 503.625
$ws.Range("H132").Value = 4998.3335
$ws.Range("I132").Value = 4998
$ws.Range("J132").Value = 4998.5
$ws.Range("K132").Value = 14994
$ws.Range("L132").Value = 14995.5
$ws.Range("M132").Value = -12464
$ws.Range("N132").Value = -20055.5
$ws.Range("H136").Value = 5837.077
$ws.Range("I136").Value = 5654.222
$ws.Range("K136").Value = 16962.666
$ws.Range("M136").Value = -14412.666
